{"js": "// Load the body paragraphs so we can inspect their text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldTail = \"- I'm Rudra, a student from Toronto, Canada. ... I am studying computer science and play sports, too. How cool is that! ...\";\nconst newTail = \"Rudra Cantaria is studying computer science at the University of Waterloo. They love playing sports and go to the gym. They are looking for a conversation partner. They want to write a blog about their passion for computers, and share tips or ideas with people.\";\n\nconst toDelete = [];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n\n  if (text.indexOf(oldTail) !== -1) {\n    // Replace the transcript paragraph's trailing summary sentence with the\n    // new bio-style summary, keeping everything before it untouched.\n    const updated = text.split(oldTail).join(newTail);\n    para.insertText(updated, \"Replace\");\n  } else if (\n    text.indexOf(\"I love going to the gym\") !== -1 ||\n    text.indexOf(\"I also enjoy coding and creating\") !== -1 ||\n    text.indexOf(\"I'm an big fan of continuous learning\") !== -1\n  ) {\n    // These three bullet paragraphs are removed entirely.\n    toDelete.push(para);\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the trailing quick-summary sentence at the end of the transcript\n# paragraph with a new third-person bio-style summary. Everything before\n# \"There you go.\" is left untouched.\n$oldTail = \"- I'm Rudra, a student from Toronto, Canada. ... I am studying computer science and play sports, too. How cool is that! ...\"\n$newTail = \"Rudra Cantaria is studying computer science at the University of Waterloo. They love playing sports and go to the gym. They are looking for a conversation partner. They want to write a blog about their passion for computers, and share tips or ideas with people.\"\n\n$find = $d.Content.Find\n$find.Text = $oldTail\n$find.Execute($oldTail, $false, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)\n\n# Remove the three bulleted \"about me\" paragraphs entirely.\n$targets = @(\n    \"I love going to the gym, and I am always up for trying new workouts.\",\n    \"I also enjoy coding and creating, and have a strong passion for solving complex problems.\",\n    \"I'm an big fan of continuous learning and staying in touch with the latest trends/technologies in the tech world.\"\n)\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs($i)\n    $t = $p.Range.Text\n    foreach ($target in $targets) {\n        if ($t.Contains($target)) {\n            $p.Range.Delete()\n        }\n    }\n}\n"}
